$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for rows 2..416.
# The commit bumps every value in that range from 45203 (2023-10-04)
# to 45204 (2023-10-05).
for ($r = 2; $r -le 416; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45203) {
        $cell.Value2 = 45204
    }
}
